# This sheet holds a rolling log of weekly price reports for
# "Comercializadora del Agro de Limarí - Pepino dulce". Each week's report
# is a small block of 3-4 rows (quality grades Especial/Primera/Segunda/
# Tercera) inserted just above the previous weeks, pushing the whole table
# down. This commit ("Fruta / hortaliza, semanal") adds the new weekly
# report for 2022-01-27 (Excel serial date 44588).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at the top of the data block (row 246), shifting
# all the existing report rows (246:323) down to (249:326). This alone
# reproduces every "D/I/J/K/L/M/P shifted down by 3 rows" change in the
# diff for the pre-existing data, since the rest of the row content is
# identical for every record in this sheet.
$ws.Range("246:248").Insert()

# Fill in the 3 new rows with this week's reported prices.
function Set-Row($row, $d, $i, $j, $k, $l, $m, $p) {
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j
    $ws.Cells.Item($row, 11).Value = $k
    $ws.Cells.Item($row, 12).Value = $l
    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $p
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

# D=fecha(44588=2022-01-27), I=calidad, J=volumen, K=min, L=max, M=promedio, P=precio $/Kg
Set-Row 246 44588 "Primera" 500 13000 14000 13500 750
Set-Row 247 44588 "Segunda" 500 11000 12000 11500 639
Set-Row 248 44588 "Tercera" 300 8000 9000 8500 472
